$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("partidas")

# Row 14: fill RF (E14), Gols_Time1 (G14), Gols_Time2 (H14), and update Status (I14)
$ws.Range("E14").Value = "3x2"
$ws.Range("G14").Value = 3
$ws.Range("H14").Value = 2
$ws.Range("I14").Value = "Finalizado"

# Row 15: fill RF (E15), Gols_Time1 (G15), Gols_Time2 (H15), and update Status (I15)
$ws.Range("E15").Value = "3x3"
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 3
$ws.Range("I15").Value = "Finalizado"

# Update the active selection to match the saved view state
$ws.Range("L15").Select()
